# Add date support: a "Dato" (date) label in A4 next to the existing date
# value in B4, plus a new header "Cell 2 C" alongside the existing "Cell 2 B"
# header, and make "Sheet number one" the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New cells: C2 gets a third header label, A4 gets the "Dato" label that
# goes with the existing date value already sitting in B4.
$ws1.Cells.Item(2, 3).Value = "Cell 2 C"
$ws1.Cells.Item(4, 1).Value = "Dato"

# Make "Sheet number one" the active sheet and put the selection on A5.
$ws1.Activate()
$null = $ws1.Range("A5").Select()
